$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": extend table to a new row 5, update rows 2-5 ---
$ws1 = $wb.Worksheets.Item("Schedule")

# Extend formatting down to the new row 5 by copying the date/time style of
# row 4 (cols A:B) before the values are changed.
$ws1.Range("A4:B4").Copy()
$ws1.Range("A5:B5").PasteSpecial(-4122)

$ws1.Cells.Item(2, 1).Value = 46039.02083333334
$ws1.Cells.Item(2, 2).Value = 46039.1875
$ws1.Cells.Item(2, 3).Value = 4
$ws1.Cells.Item(2, 4).Value = 15.12
$ws1.Cells.Item(2, 5).Value = 347.06108775
$ws1.Cells.Item(2, 6).Value = 22.95377564484127

$ws1.Cells.Item(3, 1).Value = 46039.29166666666
$ws1.Cells.Item(3, 2).Value = 46039.5
$ws1.Cells.Item(3, 3).Value = 5
$ws1.Cells.Item(3, 4).Value = 18.9
$ws1.Cells.Item(3, 5).Value = 313.99407
$ws1.Cells.Item(3, 6).Value = 16.61344285714286

$ws1.Cells.Item(4, 1).Value = 46039.52083333334
$ws1.Cells.Item(4, 2).Value = 46039.8125
$ws1.Cells.Item(4, 3).Value = 7
$ws1.Cells.Item(4, 4).Value = 26.46
$ws1.Cells.Item(4, 5).Value = 301.5587835000001
$ws1.Cells.Item(4, 6).Value = 11.39677942176871

$ws1.Cells.Item(5, 1).Value = 46040.29166666666
$ws1.Cells.Item(5, 2).Value = 46040.79166666666
$ws1.Cells.Item(5, 3).Value = 12
$ws1.Cells.Item(5, 4).Value = 45.36
$ws1.Cells.Item(5, 5).Value = 149.74963575
$ws1.Cells.Item(5, 6).Value = 3.301358812830688


# --- Sheet "Detailed": update individual cell values/status flags ---
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Cells.Item(3, 5).Value = "ON"
$ws2.Cells.Item(4, 5).Value = "ON"
$ws2.Cells.Item(5, 5).Value = "ON"
$ws2.Cells.Item(6, 5).Value = "ON"
$ws2.Cells.Item(7, 5).Value = "ON"
$ws2.Cells.Item(8, 5).Value = "ON"
$ws2.Cells.Item(9, 5).Value = "ON"
$ws2.Cells.Item(10, 5).Value = "ON"
$ws2.Cells.Item(26, 5).Value = "OFF"
$ws2.Cells.Item(39, 2).Value = -2.77661
$ws2.Cells.Item(40, 2).Value = 14.77494
$ws2.Cells.Item(41, 3).Value = "historical"
$ws2.Cells.Item(42, 3).Value = "historical"
$ws2.Cells.Item(43, 2).Value = 56.98
$ws2.Cells.Item(43, 5).Value = "OFF"
$ws2.Cells.Item(44, 2).Value = 56.98
$ws2.Cells.Item(44, 5).Value = "OFF"
$ws2.Cells.Item(45, 2).Value = 46.39404
$ws2.Cells.Item(45, 5).Value = "OFF"
$ws2.Cells.Item(46, 2).Value = 36.2
$ws2.Cells.Item(46, 5).Value = "OFF"
$ws2.Cells.Item(47, 2).Value = 47.10776
$ws2.Cells.Item(47, 5).Value = "OFF"
$ws2.Cells.Item(48, 2).Value = 56.98
$ws2.Cells.Item(48, 5).Value = "OFF"
$ws2.Cells.Item(49, 2).Value = 40.54
$ws2.Cells.Item(49, 5).Value = "OFF"
$ws2.Cells.Item(50, 5).Value = "OFF"
$ws2.Cells.Item(52, 2).Value = 31.24649
$ws2.Cells.Item(57, 2).Value = 56.98
$ws2.Cells.Item(58, 2).Value = 56.98
$ws2.Cells.Item(59, 2).Value = 36.2
$ws2.Cells.Item(61, 2).Value = 57.06003
$ws2.Cells.Item(62, 2).Value = 57.06003
$ws2.Cells.Item(64, 2).Value = 24.28595
$ws2.Cells.Item(64, 5).Value = "ON"
$ws2.Cells.Item(65, 2).Value = 28.49947
$ws2.Cells.Item(66, 2).Value = 0.51
$ws2.Cells.Item(67, 2).Value = 1.82692
$ws2.Cells.Item(68, 2).Value = 0.7
$ws2.Cells.Item(69, 2).Value = 0.7
$ws2.Cells.Item(70, 2).Value = 35.88
$ws2.Cells.Item(71, 2).Value = 35.88
$ws2.Cells.Item(72, 2).Value = 26.63005
$ws2.Cells.Item(73, 2).Value = 22.07
$ws2.Cells.Item(74, 2).Value = 0.7
$ws2.Cells.Item(75, 2).Value = 0.65374
$ws2.Cells.Item(76, 2).Value = 0.51
$ws2.Cells.Item(77, 2).Value = 0
$ws2.Cells.Item(78, 2).Value = -4.42723
$ws2.Cells.Item(79, 2).Value = -0.83768
$ws2.Cells.Item(80, 2).Value = 0.00025
$ws2.Cells.Item(81, 2).Value = -0.32843
$ws2.Cells.Item(82, 2).Value = -1.79043
$ws2.Cells.Item(83, 2).Value = -6.5608
$ws2.Cells.Item(84, 2).Value = -6.56595
$ws2.Cells.Item(85, 2).Value = -3.90305
$ws2.Cells.Item(86, 2).Value = -0.84376
$ws2.Cells.Item(87, 2).Value = 0.00032
$ws2.Cells.Item(89, 2).Value = 46.19147
$ws2.Cells.Item(90, 2).Value = 46.26191
$ws2.Cells.Item(91, 2).Value = 55.49501
$ws2.Cells.Item(92, 2).Value = 46.85269
$ws2.Cells.Item(94, 2).Value = 55.27992
